$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace host entries with a single network entry
$ws.Range("A2").Value = "Net_10.1.0.0_24"
$ws.Range("A3").Value = $null
$ws.Range("A4").Value = $null
$ws.Range("A5").Value = $null

# Update the active selection to E8 as recorded by the workbook view
$ws.Range("E8").Select()
